$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 99: new timelog entry on 2015-04-26 (serial 42120), 11:00 - 15:00
$ws.Range("A99").Value = 42120
$ws.Range("B99").Value = 0.45833333333333331
$ws.Range("C99").Value = 0.625
$ws.Range("B99:C99").NumberFormat = "h:mm"
$ws.Range("E99").Value = "thesis chapter fundamentals, terms"

# Row 100: new timelog entry on 2015-04-26 (serial 42120), 16:00 - 19:30
$ws.Range("A100").Value = 42120
$ws.Range("B100").Value = 0.66666666666666663
$ws.Range("C100").Value = 0.8125
$ws.Range("B100:C100").NumberFormat = "h:mm"
$ws.Range("E100").Value = "thesis chapter fundamentals, terms"

# Update the active selection to reflect where the author left the cursor
$ws.Range("E111").Select()
